# Refresh cryptocurrency Price / Volume(1h) figures (and two swapped rank rows)
# to match the "Updated cryptos list" GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.642.20'
$ws.Cells.Item(2, 5).Value = '  +1.05%  '
$ws.Cells.Item(3, 4).Value = '1.852.88'
$ws.Cells.Item(3, 5).Value = '  +0.39%  '
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '0.9984'
$ws.Cells.Item(4, 5).Value = '  +0.03%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '240.86'
$ws.Cells.Item(5, 5).Value = '  +0.16%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.6313'
$ws.Cells.Item(6, 5).Value = '  +0.59%  '
$ws.Cells.Item(7, 5).Value = '  -0.01%  '
$ws.Cells.Item(8, 5).Value = '  -1.17%  '
$ws.Cells.Item(9, 5).Value = '  +0.36%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '25.10'
$ws.Cells.Item(10, 5).Value = '  +2.63%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.07750'
$ws.Cells.Item(11, 5).Value = '  -0.01%  '
$ws.Cells.Item(12, 4).Value = '1.853.47'
$ws.Cells.Item(12, 5).Value = '  +0.41%  '
$ws.Cells.Item(13, 5).Value = '  +0.68%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '0.6830'
$ws.Cells.Item(14, 5).Value = '  +0.80%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.00001029'
$ws.Cells.Item(15, 5).Value = '  -0.78%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '82.90'
$ws.Cells.Item(16, 5).Value = '  -0.15%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '6.335'
$ws.Cells.Item(17, 5).Value = '  +3.91%  '
$ws.Cells.Item(18, 4).Value = '29.617.33'
$ws.Cells.Item(18, 5).Value = '  +1.01%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '230.69'
$ws.Cells.Item(19, 5).Value = '  +0.78%  '
$ws.Cells.Item(20, 5).Value = '  +0.76%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '0.9998'
$ws.Cells.Item(21, 5).Value = '  +0.04%  '
$ws.Cells.Item(22, 5).Value = '  +1.67%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '0.9989'
$ws.Cells.Item(23, 5).Value = '  -0.32%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '159.52'
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '8.526'
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '17.59'
$ws.Cells.Item(27, 5).Value = '  -0.33%  '
$ws.Cells.Item(28, 2).Value = 'Hedera'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '0.06583'
$ws.Cells.Item(28, 5).Value = '  +15.80%  '
$ws.Cells.Item(29, 2).Value = 'Toncoin'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '1.467'
$ws.Cells.Item(29, 5).Value = '  +3.69%  '
$ws.Cells.Item(30, 5).Value = '  +1.13%  '
$ws.Cells.Item(31, 5).Value = '  +0.28%  '
$ws.Cells.Item(32, 5).Value = '  +1.87%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '1.854'
$ws.Cells.Item(33, 5).Value = '  +1.96%  '
$ws.Cells.Item(34, 5).Value = '  -0.53%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.7007'
$ws.Cells.Item(35, 5).Value = '  +0.63%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '2.566'
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '0.01871'
$ws.Cells.Item(38, 2).Value = 'Maker'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(38, 4).Value = '1.260.38'
$ws.Cells.Item(38, 5).Value = '  +1.90%  '
$ws.Cells.Item(39, 2).Value = 'MXToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '2.841'
$ws.Cells.Item(39, 5).Value = '  +4.56%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '6.761'
$ws.Cells.Item(40, 5).Value = '  +5.44%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.9426'
$ws.Cells.Item(41, 5).Value = '  +4.73%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '1.001'
$ws.Cells.Item(42, 5).Value = '  +0.19%  '
$ws.Cells.Item(43, 4).Value = '2.002.42'
$ws.Cells.Item(43, 5).Value = '  -0.10%  '
$ws.Cells.Item(44, 5).Value = '  +0.07%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '66.30'
$ws.Cells.Item(45, 5).Value = '  +1.30%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '1.741'
$ws.Cells.Item(46, 5).Value = '  +4.17%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '7.126'
$ws.Cells.Item(47, 5).Value = '  +0.06%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '0.00000000117'
$ws.Cells.Item(48, 5).Value = '  +1.68%  '
$ws.Cells.Item(49, 5).Value = '  +1.39%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '9.031'
$ws.Cells.Item(50, 5).Value = '  +0.29%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.3958'
$ws.Cells.Item(51, 5).Value = '  -0.86%  '
